# Apply the commit's changes:
#  - Database views referenced in the "Tables" sheet must be schema-qualified
#    (prefix every "view_component_*" value in column B with "reference.").
#  - The "Tables" sheet becomes the active/selected sheet (instead of "ROOT"),
#    with its view reset to the top-left cell and selection on B2.
#  - Column B on "Tables" is widened (no longer auto bestFit) to comfortably
#    fit the longer, schema-qualified view names.

$wb = $excel.ActiveWorkbook

$rootWs   = $wb.Worksheets.Item("ROOT")
$tablesWs = $wb.Worksheets.Item("Tables")

# 1. Prefix every view_component_* reference in column B with "reference."
for ($row = 2; $row -le 31; $row++) {
    $cell = $tablesWs.Cells.Item($row, 2)
    $current = $cell.Value
    if ($current -ne $null -and $current.ToString().StartsWith("view_component_")) {
        $cell.Value = "reference." + $current
    }
}

# 2. Widen column B to fit the new, longer values (fixed width, not bestFit).
$tablesWs.Columns.Item(2).ColumnWidth = 48.36328125

# 3. Move the active tab / selection from ROOT to Tables.
$tablesWs.Select()
$tablesWs.Range("A1").Select()
$tablesWs.Range("B2").Select()

$wb.Save()
